$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value (never auto-converted to a date/number
# by Excel's input parser) into a cell, while preserving the numeric-style
# formatting ($styleSourceRef's cell style) the target cell should end up
# with.
function Set-TextValue($ws, $cellRef, $text, $styleSourceRef) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
    $ws.Range($styleSourceRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# 1) "Bad Drivers" table: insert a new data row (row 4) for a newly observed
#    driver, pushing the existing "Totals:" row (and everything below it)
#    down by one.
# ---------------------------------------------------------------------------
$ws.Range("A4").EntireRow.Insert()

# Copy formatting from the row above (the existing data row) onto the new row
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 98.90000000000001

# Update the figures that changed on the existing first data row (row 3)
$ws.Range("C3").Value = 668
$ws.Range("D3").Value = 97.5

# Update the "Totals:" row (now row 5)
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 672

# ---------------------------------------------------------------------------
# 2) "Good Drivers" table: insert two new data rows right after the column
#    header row (now row 12) for two newly observed drivers, then append
#    five more new rows at the bottom of the table.
# ---------------------------------------------------------------------------
$ws.Range("A13:A14").EntireRow.Insert()

# Copy formatting from the data row that follows (the previously-first data
# row, now pushed down) onto the two freshly inserted rows.
$ws.Range("A15:E15").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Range("A14:E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B13").Value = 11128
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = ""

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B14").Value = 486214
$ws.Range("D14").Value = 99.90000000000001
Set-TextValue $ws "E14" "2024-11-10" "E15"

# Update the sample counts for the (now shifted) pre-existing driver rows
$ws.Range("B15").Value = 11140
$ws.Range("B16").Value = 14487

# Append five brand-new driver rows after the last existing data row (now
# row 17), re-using its formatting.
$ws.Range("A17:E17").Copy()
$ws.Range("A18:E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B18").Value = 79953
$ws.Range("D18").Value = 99.90000000000001
Set-TextValue $ws "E18" "2021-08-18" "E17"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B19").Value = 35355
$ws.Range("D19").Value = 100
Set-TextValue $ws "E19" "2021-04-27" "E17"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B20").Value = 65425
$ws.Range("D20").Value = 100
Set-TextValue $ws "E20" "2020-08-05" "E17"

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B21").Value = 117653
$ws.Range("D21").Value = 100
Set-TextValue $ws "E21" "2020-01-06" "E17"

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B22").Value = 56018
$ws.Range("D22").Value = 100
Set-TextValue $ws "E22" "2019-12-14" "E17"

# ---------------------------------------------------------------------------
# 3) The report template always reserves a trailing band of blank rows below
#    each table, and spacer columns F:J. Touch the bottom-right corner cell
#    (a cheap no-op format write) so the sheet's used range / dimension
#    extends all the way out to J27, matching the refreshed template.
# ---------------------------------------------------------------------------
$ws.Range("J27").Font.Bold = $false
